# Update countries & provincias Spain
# Applies the shared-string reorder / data refresh described in the commit
# by writing the resulting cell values directly onto the "Pais" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp update (row 1) -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 08:35"

# --- Data rows that changed -------------------------------------------------
# Each entry: Row, Country (A), Casos totales (B), Nuevos casos (C),
# Casos activos (D), Recuperados (E), Casos criticos (F), Muertes hoy (G),
# Muertes (H)
$rows = @(
    @{ Row = 39;  A = "Ucrania";                        B = 21584; C = 339; D = 7575;  E = 13365; F = 0; G = 21; H = 644 },
    @{ Row = 48;  A = "Afganistan";                      B = 11831; C = 658; D = 1128;  E = 10483; F = 0; G = 1;  H = 220 },
    @{ Row = 49;  A = "Dinamarca";                       B = 11387; C = 0;   D = 9964;  E = 860;   F = 0; G = 0;  H = 563 },
    @{ Row = 50;  A = "Corea del Sur";                   B = 11225; C = 19;  D = 10275; E = 681;   F = 0; G = 2;  H = 269 },
    @{ Row = 51;  A = "Serbia";                          B = 11193; C = 0;   D = 5920;  E = 5034;  F = 0; G = 0;  H = 239 },
    @{ Row = 52;  A = "Panama";                          B = 11183; C = 0;   D = 6279;  E = 4594;  F = 0; G = 0;  H = 310 },
    @{ Row = 197; A = "Fiyi";                            B = 18;    C = 0;   D = 15;    E = 3;     F = 0; G = 0;  H = 0   },
    @{ Row = 198; A = "Curazao";                         B = 18;    C = 0;   D = 14;    E = 3;     F = 0; G = 0;  H = 1   },
    @{ Row = 199; A = "Santa Lucia";                     B = 18;    C = 0;   D = 18;    E = 0;     F = 0; G = 0;  H = 0   },
    @{ Row = 201; A = "Nueva Caledonia";                 B = 18;    C = 0;   D = 18;    E = 0;     F = 0; G = 0;  H = 0   },
    @{ Row = 215; A = "San Bartolome";                   B = 6;     C = 0;   D = 6;     E = 0;     F = 0; G = 0;  H = 0   },
    @{ Row = 216; A = "Bonaire, San Eustaquio y Saba";   B = 6;     C = 0;   D = 6;     E = 0;     F = 0; G = 0;  H = 0   }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}
